$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.403.89"
$ws.Range("E2").Value = "  +0.15%  "

# Row 3
$ws.Range("D3").Value = "'1.849.53"
$ws.Range("E3").Value = "  +0.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'240.65"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").Value = "'0.6301"
$ws.Range("E6").Value = "  -0.18%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.07675"
$ws.Range("E8").Value = "  +1.86%  "

# Row 9
$ws.Range("D9").Value = "'0.2940"
$ws.Range("E9").Value = "  -0.61%  "

# Row 10
$ws.Range("D10").Value = "'24.50"
$ws.Range("E10").Value = "  +0.27%  "

# Row 11
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("D12").Value = "'1.845.34"
$ws.Range("E12").Value = "  -0.31%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.018"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.00001092"
$ws.Range("E14").Value = "  +9.05%  "

# Row 15
$ws.Range("D15").Value = "'0.6801"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("D16").Value = "'83.59"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").Value = "'2.097.38"
$ws.Range("E17").Value = "  -7.37%  "

# Row 18
$ws.Range("D18").Value = "'6.148"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("D19").Value = "'29.414.13"

# Row 20
$ws.Range("D20").Value = "'229.61"
$ws.Range("E20").Value = "  +0.99%  "

# Row 21
$ws.Range("D21").Value = "'12.46"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").Value = "'7.444"
$ws.Range("E23").Value = "  -1.18%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").Value = "'157.31"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("E26").Value = "  -0.44%  "

# Row 27
$ws.Range("D27").Value = "'8.363"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").Value = "'17.68"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.467"
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.311"
$ws.Range("E30").Value = "  +4.43%  "

# Row 31
$ws.Range("D31").Value = "'0.05676"
$ws.Range("E31").Value = "  +0.02%  "

# Row 32
$ws.Range("D32").Value = "'4.114"
$ws.Range("E32").Value = "  -0.29%  "

# Row 33
$ws.Range("D33").Value = "'4.049"
$ws.Range("E33").Value = "  +0.89%  "

# Row 34
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").Value = "'0.7113"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("E37").Value = "  -0.20%  "

# Row 38
$ws.Range("E38").Value = "  -0.11%  "

# Row 39
$ws.Range("D39").Value = "'1.232.44"
$ws.Range("E39").Value = "  -1.79%  "

# Row 40
$ws.Range("E40").Value = "  -0.88%  "

# Row 41
$ws.Range("D41").Value = "'6.473"
$ws.Range("E41").Value = "  +4.29%  "

# Row 42
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'2.006.50"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.49"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.22"
$ws.Range("E46").Value = "  -0.28%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000121"
$ws.Range("E47").Value = "  +3.87%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.161"
$ws.Range("E48").Value = "  +1.75%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4015"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.024"
$ws.Range("E50").Value = "  -0.99%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.689"
$ws.Range("E51").Value = "  +0.04%  "

